# #17 - Validacao de ramo de atividade
#
# The "Ramo de Atividade do Imovel" sheet gets a new "Posição" column
# (matching the same column already present on the Cliente / Imovel /
# Servicos / Medidor / Anormalidade Imovel sheets), and the active
# workbook tab/selection state moves onto that sheet.

$wb = $excel.ActiveWorkbook

$shRamo     = $wb.Worksheets.Item("Ramo de Atividade do Imovel")
$shCliente  = $wb.Worksheets.Item("Cliente")
$shImovel   = $wb.Worksheets.Item("Imovel")

# ---------------------------------------------------------------
# 1) Ramo de Atividade do Imovel: insert a "Posição" column at D,
#    pushing the old D (Padrão) / E (Comentário) columns to E / F.
# ---------------------------------------------------------------
$shRamo.Columns.Item(4).Insert() | Out-Null

# Match column C's width for the freshly inserted column D.
$shRamo.Columns.Item(4).ColumnWidth = 7.5

# Header + data for the new column.
$shRamo.Range("D2").Value = "Posição"
$shRamo.Range("D3").Value = 1
$shRamo.Range("D4").Formula = "=SUM(C3,D3)"
$shRamo.Range("D5").Formula = "=SUM(C4,D4)"

# Pick up the same "Posição" formula-cell formatting used elsewhere
# in the workbook (centered, themed border) instead of the plain
# copy-from-the-left formatting that Insert() applied.
$shCliente.Range("D4").Copy() | Out-Null
$shRamo.Range("D4").PasteSpecial(-4122) | Out-Null
$shCliente.Range("D5").Copy() | Out-Null
$shRamo.Range("D5").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------
# 2) Move the live selection around on the other sheets touched by
#    this edit, then land on Ramo de Atividade do Imovel so it
#    becomes the workbook's active tab.
# ---------------------------------------------------------------
$shCliente.Activate() | Out-Null
$shCliente.Range("B8").Select() | Out-Null

$shImovel.Activate() | Out-Null
$shImovel.Range("B27").Select() | Out-Null

$shRamo.Activate() | Out-Null
$shRamo.Range("D4:D5").Select() | Out-Null
